## Jan 14, 2020 Update
## - Widen the Weibull-mixture shape parameter range used in the
##   "Score Function Simulations" sheet (shape=c(1.5,3) -> shape=c(1.1,14))
## - Bump the iterations (N) parameter on the same sheet from 2 to 50

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Score Function Simulations")
$ws.Activate()

# Update the "parameters" description string for the rweibullmix distribution
$ws.Range("B3").Value = "list(lambda=c(0.25,0.75), shape=c(1.1,14), scale=c(0.1,5))"

# Update N from 2 to 50
$ws.Range("B6").Value = 50

# Leave the selection on the cell that was last edited
[void]$ws.Range("B6").Select()
